$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.353.24'
$ws.Range("E2").Value = '  -1.38%  '
# Row 3
$ws.Range("D3").Value = '1.592.19'
$ws.Range("E3").Value = '  -0.45%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.57%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '
# Row 6
$ws.Range("E6").Value = '  -1.99%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.56%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0612'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.10%  '
# Row 9
$ws.Range("E9").Value = '  -0.46%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.16%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.57%  '
# Row 12
$ws.Range("E12").Value = '  -0.50%  '
# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.29%  '
# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.581.80'
$ws.Range("E14").Value = '  -0.82%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.518'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.36%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.47%  '
# Row 17
$ws.Range("D17").Value = '26.350.63'
$ws.Range("E17").Value = '  -1.29%  '
# Row 18
$ws.Range("E18").Value = '  -1.34%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.73%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '212.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.65%  '
# Row 21
$ws.Range("E21").Value = '  -0.60%  '
# Row 22
$ws.Range("E22").Value = '  -0.19%  '
# Row 23
$ws.Range("E23").Value = '  -1.97%  '
# Row 24
$ws.Range("E24").Value = '  -1.19%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.60%  '
# Row 26
$ws.Range("E26").Value = '  -0.52%  '
# Row 27
$ws.Range("E27").Value = '  -1.08%  '
# Row 28
$ws.Range("E28").Value = '  -1.37%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.28%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0505'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.31%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.75%  '
# Row 32
$ws.Range("E32").Value = '  -0.84%  '
# Row 33
$ws.Range("E33").Value = '  +0.79%  '
# Row 34
$ws.Range("D34").Value = '1.297.48'
$ws.Range("E34").Value = '  +1.72%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.612'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.96%  '
# Row 36
$ws.Range("E36").Value = '  -1.79%  '
# Row 37
$ws.Range("E37").Value = '  -1.02%  '
# Row 38
$ws.Range("E38").Value = '  -0.83%  '
# Row 39
$ws.Range("E39").Value = '  -10.13%  '
# Row 40
$ws.Range("E40").Value = '  -1.26%  '
# Row 41
$ws.Range("E41").Value = '  -0.60%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.16%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.21%  '
# Row 44
$ws.Range("E44").Value = '  -2.53%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.761'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.07%  '
# Row 46
$ws.Range("D46").Value = '1.727.00'
$ws.Range("E46").Value = '  -0.43%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.20%  '
# Row 48
$ws.Range("E48").Value = '  -3.42%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.98%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.34%  '
# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.01%  '
